$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a literal text value into a cell without Excel's
# "looks like a number" auto-conversion kicking in (this matters because
# several of the new values start with "+" followed by digits, which Excel
# would otherwise coerce to a numeric value). We build the text through a
# formula that evaluates to a string, then "bake" it into a static value via
# Copy / PasteSpecial(values-only) so the final cell is a plain text cell,
# identical to typing it in and keeping the existing cell style untouched.
function Set-TextValue {
    param($address, $text)
    $helper = $ws.Range("ZZ1")
    $helper.Formula = '="' + $text + '"'
    $helper.Copy()
    $ws.Range($address).PasteSpecial(-4163)  # xlPasteValues
    $helper.ClearContents()
}

# Row 4: the garbled phone number is actually the same number as rows 2-3,
# so just reuse that exact (already-correct) text instead of retyping it.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4163)  # xlPasteValues

# Rows 5-7: correct the invalid-looking phone numbers, and reset their
# (clearly bogus / placeholder) amount columns to small sanity-check values.
Set-TextValue "A5" "+212611250473"
$ws.Range("B5").Value = 123
$ws.Range("C5").Value = 123

Set-TextValue "A6" "+21261125047"
$ws.Range("B6").Value = 234
$ws.Range("C6").Value = 234

Set-TextValue "A7" "+2122323"
$ws.Range("B7").Value = 234
$ws.Range("C7").Value = 234

# Row 8 (previously blank): add one more corrected phone number entry, with
# no amounts yet.
Set-TextValue "A8" "+2126111111111"

# Column A needs to be noticeably wider now that it holds longer numbers.
$ws.Columns("A").ColumnWidth = 21.8

# Leave the selection on the newly added last row, like a user would after
# finishing data entry there.
$ws.Range("A8").Select()
